$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 37 (pushes existing rows 37-88 down to 38-89)
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A37").Value = 4
$ws.Range("B37").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C37").Value = "Los Lagos"
$ws.Range("D37").Value = 44757
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 100112026
$ws.Range("G37").Value = "Haba"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 80
$ws.Range("K37").Value = 22000
$ws.Range("L37").Value = 22000
$ws.Range("M37").Value = 22000
$ws.Range("N37").Value = "`$/saco 25 kilos"
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 880
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
